$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A36").Value = 45671
$ws.Range("A36").NumberFormat = $ws.Range("A35").NumberFormat
$ws.Range("B36").Value = "fleshing out all made areas so far"
$ws.Range("C36").Value = 6

$ws.Range("B37").Select()
